$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 72.266001
$ws.Range("H2").Value = 216.798003
$ws.Range("I2").Value = 0.2949652269937106
$ws.Range("J2").Value = 0.2949652269937106
$ws.Range("M2").Value = 17.46627766666667
$ws.Range("N2").Value = 52.398833
$ws.Range("O2").Value = 0.2609791297364465
$ws.Range("P2").Value = 0.2609791297364465
$ws.Range("Q2").Value = 1262.218039325611
$ws.Range("R2").Value = 11359.9623539305
$ws.Range("S2").Value = 0.07697976824333198
$ws.Range("T2").Value = 0.07697976824333198
$ws.Range("G3").Value = 72.266001
$ws.Range("H3").Value = 216.798003
$ws.Range("I3").Value = 0.2949652269937106
$ws.Range("J3").Value = 0.2949652269937106
$ws.Range("O3").Value = 0.03718296798122674
$ws.Range("P3").Value = 0.03718296798122674
$ws.Range("Q3").Value = 179.834352995839
$ws.Range("R3").Value = 1618.509176962551
$ws.Range("S3").Value = 0.01096768259088242
$ws.Range("T3").Value = 0.01096768259088242
$ws.Range("G4").Value = 72.266001
$ws.Range("H4").Value = 216.798003
$ws.Range("I4").Value = 0.2949652269937106
$ws.Range("J4").Value = 0.2949652269937106
$ws.Range("M4").Value = 46.97117233333334
$ws.Range("N4").Value = 140.913517
$ws.Range("O4").Value = 0.7018379022823268
$ws.Range("P4").Value = 0.7018379022823268
$ws.Range("Q4").Value = 3394.418786811839
$ws.Range("R4").Value = 30549.76908130655
$ws.Range("S4").Value = 0.2070177761594962
$ws.Range("T4").Value = 0.2070177761594962
$ws.Range("H5").Value = 410.023338
$ws.Range("I5").Value = 0.5578585839920717
$ws.Range("J5").Value = 0.5578585839920718
$ws.Range("M5").Value = 17.46627766666667
$ws.Range("N5").Value = 52.398833
$ws.Range("O5").Value = 0.2609791297364465
$ws.Range("P5").Value = 0.2609791297364465
$ws.Range("Q5").Value = 2387.193823773839
$ws.Range("R5").Value = 21484.74441396455
$ws.Range("S5").Value = 0.1455894477662572
$ws.Range("T5").Value = 0.1455894477662572
$ws.Range("H6").Value = 410.023338
$ws.Range("I6").Value = 0.5578585839920717
$ws.Range("J6").Value = 0.5578585839920718
$ws.Range("O6").Value = 0.03718296798122674
$ws.Range("P6").Value = 0.03718296798122674
$ws.Range("Q6").Value = 340.1151333595273
$ws.Range("S6").Value = 0.02074283786662969
$ws.Range("T6").Value = 0.02074283786662969
$ws.Range("H7").Value = 410.023338
$ws.Range("I7").Value = 0.5578585839920717
$ws.Range("J7").Value = 0.5578585839920718
$ws.Range("M7").Value = 46.97117233333334
$ws.Range("N7").Value = 140.913517
$ws.Range("O7").Value = 0.7018379022823268
$ws.Range("P7").Value = 0.7018379022823268
$ws.Range("Q7").Value = 6419.758956628861
$ws.Range("R7").Value = 57777.83060965974
$ws.Range("S7").Value = 0.3915262983591848
$ws.Range("T7").Value = 0.3915262983591848
$ws.Range("G8").Value = 36.057927
$ws.Range("H8").Value = 108.173781
$ws.Range("I8").Value = 0.1471761890142177
$ws.Range("J8").Value = 0.1471761890142177
$ws.Range("M8").Value = 17.46627766666667
$ws.Range("N8").Value = 52.398833
$ws.Range("O8").Value = 0.2609791297364465
$ws.Range("P8").Value = 0.2609791297364465
$ws.Range("Q8").Value = 629.797765066397
$ws.Range("R8").Value = 5668.179885597573
$ws.Range("S8").Value = 0.03840991372685729
$ws.Range("T8").Value = 0.03840991372685729
$ws.Range("G9").Value = 36.057927
$ws.Range("H9").Value = 108.173781
$ws.Range("I9").Value = 0.1471761890142177
$ws.Range("J9").Value = 0.1471761890142177
$ws.Range("O9").Value = 0.03718296798122674
$ws.Range("P9").Value = 0.03718296798122674
$ws.Range("Q9").Value = 89.730355667753
$ws.Range("R9").Value = 807.573201009777
$ws.Range("S9").Value = 0.005472447523714631
$ws.Range("T9").Value = 0.005472447523714631
$ws.Range("G10").Value = 36.057927
$ws.Range("H10").Value = 108.173781
$ws.Range("I10").Value = 0.1471761890142177
$ws.Range("J10").Value = 0.1471761890142177
$ws.Range("M10").Value = 46.97117233333334
$ws.Range("N10").Value = 140.913517
$ws.Range("O10").Value = 0.7018379022823268
$ws.Range("P10").Value = 0.7018379022823268
$ws.Range("Q10").Value = 1693.683103099753
$ws.Range("R10").Value = 15243.14792789778
$ws.Range("S10").Value = 0.1032938277636458
$ws.Range("T10").Value = 0.1032938277636458
